# Rename the embedded logo pictures:
#   - Pearson logo inline pictures (in both footers) : image1.png -> image2.png
#   - BTec logo inline pictures (in both headers)     : image2.jpg -> image1.jpg
# The `name` attribute lives on <wp:docPr> and <pic:cNvPr>, which are not
# reachable through InlineShape.Name (that property is not wired to the
# underlying OOXML in this host), so we rebuild each drawing's XML via
# InlineShape.Delete() + Range.InsertXML() with the corrected name.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Replace-LogoDrawing($range, $extentCx, $extentCy, $descr, $docPrId, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $r = $shape.Range
    $shape.Delete()

    $frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
        'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
        'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
        'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
        'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
        '<w:body><w:p><w:r><w:drawing>' +
        '<wp:inline distB="0" distT="0" distL="0" distR="0">' +
        '<wp:extent cx="' + $extentCx + '" cy="' + $extentCy + '"/>' +
        '<wp:effectExtent b="0" l="0" r="0" t="0"/>' +
        '<wp:docPr descr="' + $descr + '" id="' + $docPrId + '" name="' + $newName + '"/>' +
        '<a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
        '<pic:pic><pic:nvPicPr>' +
        '<pic:cNvPr descr="' + $descr + '" id="0" name="' + $newName + '"/>' +
        '<pic:cNvPicPr preferRelativeResize="0"/>' +
        '</pic:nvPicPr>' +
        '<pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/>' +
        '<a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
        '<pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="' + $extentCx + '" cy="' + $extentCy + '"/></a:xfrm>' +
        '<a:prstGeom prst="rect"/><a:ln/></pic:spPr>' +
        '</pic:pic></a:graphicData></a:graphic>' +
        '</wp:inline></w:drawing></w:r></w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($frag)
}

$pearsonDescr = 'Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png'
$btecDescr = 'BTec_Logo-Orange'

# Footers hold the Pearson logo: image1.png -> image2.png
Replace-LogoDrawing $sec.Footers.Item(2).Range 952500 285750 $pearsonDescr 2 'image2.png'  # footer1.xml (first page)
Replace-LogoDrawing $sec.Footers.Item(1).Range 952500 285750 $pearsonDescr 4 'image2.png'  # footer2.xml (default)

# Headers hold the BTec logo: image2.jpg -> image1.jpg
Replace-LogoDrawing $sec.Headers.Item(2).Range 914400 277792 $btecDescr 1 'image1.jpg'  # header1.xml (first page)
Replace-LogoDrawing $sec.Headers.Item(1).Range 914400 277792 $btecDescr 3 'image1.jpg'  # header2.xml (default)
